$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.275.40'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '1.786.59'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '335.65'
$ws.Range("E5").Value = '  -2.85%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '0.3809'
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("D8").Value = '0.3416'
$ws.Range("E8").Value = '  -2.92%  '
$ws.Range("D9").Value = '48.36'
$ws.Range("E9").Value = '  -3.59%  '
$ws.Range("D10").Value = '1.197'
$ws.Range("E10").Value = '  -3.17%  '
$ws.Range("D11").Value = '0.07490'
$ws.Range("E11").Value = '  -3.51%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '21.91'
$ws.Range("E13").Value = '  -2.82%  '
$ws.Range("D14").Value = '6.457'
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").Value = '1.787.23'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("D16").Value = '7.079'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  -2.66%  '
$ws.Range("D18").Value = '0.06652'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").Value = '83.83'
$ws.Range("E19").Value = '  -3.46%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '6.628'
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").Value = '17.33'
$ws.Range("E22").Value = '  -2.88%  '
$ws.Range("D23").Value = '27.263.48'
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("D24").Value = '12.36'
$ws.Range("E24").Value = '  -6.33%  '
$ws.Range("D25").Value = '2.408'
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("E26").Value = '  -1.15%  '
$ws.Range("D27").Value = '2.542'
$ws.Range("E27").Value = '  -5.12%  '
$ws.Range("D28").Value = '21.30'
$ws.Range("E28").Value = '  -3.74%  '
$ws.Range("D29").Value = '153.51'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").Value = '1.989.74'
$ws.Range("D32").Value = '4.021'
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("D33").Value = '6.085'
$ws.Range("E33").Value = '  -4.74%  '
$ws.Range("D34").Value = '0.08711'
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").Value = '13.33'
$ws.Range("E35").Value = '  -4.27%  '
$ws.Range("D36").Value = '1.655'
$ws.Range("E36").Value = '  -3.82%  '
$ws.Range("D37").Value = '0.6941'
$ws.Range("E37").Value = '  -2.18%  '
$ws.Range("D38").Value = '5.446'
$ws.Range("D39").Value = '0.2208'
$ws.Range("E39").Value = '  -2.81%  '
$ws.Range("D40").Value = '0.06329'
$ws.Range("E40").Value = '  -3.25%  '
$ws.Range("D41").Value = '8.800'
$ws.Range("E41").Value = '  -2.31%  '
$ws.Range("D42").Value = '0.02343'
$ws.Range("E42").Value = '  -3.28%  '
$ws.Range("D43").Value = '1.237'
$ws.Range("E43").Value = '  -4.30%  '
$ws.Range("D44").Value = '14.42'
$ws.Range("E44").Value = '  -3.69%  '
$ws.Range("D45").Value = '0.6525'
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("D47").Value = '3.850'
$ws.Range("E47").Value = '  -3.26%  '
$ws.Range("D48").Value = '2.146'
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").Value = '129.03'
$ws.Range("D50").Value = '0.07134'
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("D51").Value = '78.98'
$ws.Range("E51").Value = '  -2.17%  '
